# This script appends 45 new test-data rows (rows 102-146) to the
# "master-reg_center_device" sheet, reproducing the pattern already present
# in the sheet (regcntr_id cycling 10002-10010, device_id incrementing from
# 3000121, lang_code "eng", is_active TRUE, cr_by "superadmin",
# cr_dtimes "now()"). It also updates the current selection to the newly
# added block and sets the print orientation, mirroring the workbook's
# saved view/page-setup state after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cycleA   = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)
$startRow = 102
$startB   = 3000121
$numRows  = 45

for ($i = 0; $i -lt $numRows; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $cycleA[$i % 9]
    $ws.Cells.Item($r, 2).Value = $startB + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Select the newly-added range, matching the saved selection/topLeftCell
# state recorded in the workbook after the edit.
$ws.Range("A102:F146").Select()

# The saved worksheet also carries an explicit portrait page setup.
$ws.PageSetup.Orientation = 1
